# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note text in A1 ---
$hoja1 = $wb.Worksheets.Item("Hoja1")

$oldText = $hoja1.Range("A1").Value2
$newText = $oldText.Replace(
    "1000 Bs = 3.27 = 12467.32 pesos",
    "1000 Bs = 3.28 = 12544.26 pesos"
).Replace(
    "12467.32 pesos = 3.26 = 964.61 Bs",
    "12544.26 pesos = 3.26 = 969.0 Bs"
)
$hoja1.Range("A1").Value = $newText

# --- tasas: update the rate cells N10, O10, N12, O12 ---
$tasas = $wb.Worksheets.Item("tasas")

$tasas.Range("N10").Value = 305
$tasas.Range("O10").Value = 3826
$tasas.Range("N12").Value = 3844.99
$tasas.Range("O12").Value = 297.011
